$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.161.42'

$ws.Range("D3").Value = '1.562.90'
$ws.Range("E3").Value = '  -0.83%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.67'
$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3790'
$ws.Range("E7").Value = '  +2.69%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3288'
$ws.Range("E8").Value = '  -1.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.74'
$ws.Range("E9").Value = '  -8.58%  '

$ws.Range("E10").Value = '  -0.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07377'
$ws.Range("E11").Value = '  -2.35%  '

$ws.Range("E12").Value = '  -0.07%  '

$ws.Range("E13").Value = '  -3.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.837'
$ws.Range("E14").Value = '  -2.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.888'
$ws.Range("E15").Value = '  -0.94%  '

$ws.Range("D16").Value = '1.565.20'
$ws.Range("E16").Value = '  -0.47%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001094'
$ws.Range("E17").Value = '  -2.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06628'
$ws.Range("E18").Value = '  -1.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.74'
$ws.Range("E19").Value = '  -2.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.469'
$ws.Range("E20").Value = '  +1.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.16'
$ws.Range("E22").Value = '  -2.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.74'
$ws.Range("E23").Value = '  -2.35%  '

$ws.Range("D24").Value = '22.161.89'
$ws.Range("E24").Value = '  -1.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.262'
$ws.Range("E25").Value = '  -5.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.552'
$ws.Range("E26").Value = '  -3.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.31'
$ws.Range("E27").Value = '  +0.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.12'
$ws.Range("E28").Value = '  -2.98%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.846'
$ws.Range("E29").Value = '  -3.12%  '

$ws.Range("D30").Value = '1.742.40'
$ws.Range("E30").Value = '  -0.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.51'
$ws.Range("E31").Value = '  -3.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.127'
$ws.Range("E32").Value = '  +3.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.057'
$ws.Range("E33").Value = '  -1.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.866'
$ws.Range("E34").Value = '  -6.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.386'
$ws.Range("E35").Value = '  -4.75%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08223'
$ws.Range("E36").Value = '  -1.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.301'
$ws.Range("E37").Value = '  -1.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02310'
$ws.Range("E38").Value = '  -6.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06225'
$ws.Range("E39").Value = '  -2.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2144'
$ws.Range("E40").Value = '  -4.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.240'
$ws.Range("E41").Value = '  -4.09%  '

$ws.Range("E42").Value = '  -3.47%  '

$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5992'
$ws.Range("E44").Value = '  -4.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.78'
$ws.Range("E45").Value = '  -2.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.760'
$ws.Range("E46").Value = '  -0.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5798'
$ws.Range("E47").Value = '  -5.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.988'
$ws.Range("E48").Value = '  -3.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.40'
$ws.Range("E49").Value = '  -3.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.174'
$ws.Range("E50").Value = '  -3.14%  '

$ws.Range("E51").Value = '  -3.38%  '
